$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").Value = "'246.52"
$ws.Range("D2").Style = $origStyle

$origStyle = $ws.Range("E2").Style
$ws.Range("E2").Value = "'0.36%"
$ws.Range("E2").Style = $origStyle

$origStyle = $ws.Range("D3").Style
$ws.Range("D3").Value = "'26.06"
$ws.Range("D3").Style = $origStyle

$origStyle = $ws.Range("E3").Style
$ws.Range("E3").Value = "'3.20%"
$ws.Range("E3").Style = $origStyle

$origStyle = $ws.Range("D4").Style
$ws.Range("D4").Value = "'5.203"
$ws.Range("D4").Style = $origStyle

$origStyle = $ws.Range("E4").Style
$ws.Range("E4").Value = "'3.10%"
$ws.Range("E4").Style = $origStyle

$origStyle = $ws.Range("E5").Style
$ws.Range("E5").Value = "'-0.19%"
$ws.Range("E5").Style = $origStyle

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").Value = "'6.479"
$ws.Range("D6").Style = $origStyle

$origStyle = $ws.Range("E6").Style
$ws.Range("E6").Value = "'-1.17%"
$ws.Range("E6").Style = $origStyle

$origStyle = $ws.Range("D7").Style
$ws.Range("D7").Value = "'0.8133"
$ws.Range("D7").Style = $origStyle

$origStyle = $ws.Range("E7").Style
$ws.Range("E7").Value = "'-0.20%"
$ws.Range("E7").Style = $origStyle

$origStyle = $ws.Range("D8").Style
$ws.Range("D8").Value = "'0.8460"
$ws.Range("D8").Style = $origStyle

$origStyle = $ws.Range("E8").Style
$ws.Range("E8").Value = "'1.31%"
$ws.Range("E8").Style = $origStyle

$ws.Range("B9").Value = "MandalaExchangeToken"

$ws.Range("C9").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.06921"
$ws.Range("D9").Style = $origStyle

$origStyle = $ws.Range("E9").Style
$ws.Range("E9").Value = "'-0.54%"
$ws.Range("E9").Style = $origStyle

$ws.Range("B10").Value = "BitrueCoin"

$ws.Range("C10").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").Value = "'0.02823"
$ws.Range("D10").Style = $origStyle

$origStyle = $ws.Range("E10").Style
$ws.Range("E10").Value = "'-0.05%"
$ws.Range("E10").Style = $origStyle

$ws.Range("B11").Value = "BitMartToken"

$ws.Range("C11").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").Value = "'0.09379"
$ws.Range("D11").Style = $origStyle

$origStyle = $ws.Range("E11").Style
$ws.Range("E11").Value = "'-0.08%"
$ws.Range("E11").Style = $origStyle

$ws.Range("B12").Value = "BitForexToken"

$ws.Range("C12").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"

$origStyle = $ws.Range("D12").Style
$ws.Range("D12").Value = "'0.001515"
$ws.Range("D12").Style = $origStyle

$origStyle = $ws.Range("E12").Style
$ws.Range("E12").Value = "'-0.92%"
$ws.Range("E12").Style = $origStyle

$ws.Range("B13").Value = "One"

$ws.Range("C13").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"

$origStyle = $ws.Range("D13").Style
$ws.Range("D13").Value = "'0.0005989"
$ws.Range("D13").Style = $origStyle

$origStyle = $ws.Range("E13").Style
$ws.Range("E13").Value = "'0.17%"
$ws.Range("E13").Style = $origStyle

$ws.Range("B14").Value = "TigerCash"

$ws.Range("C14").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").Value = "'0.006176"
$ws.Range("D14").Style = $origStyle

$origStyle = $ws.Range("E14").Style
$ws.Range("E14").Value = "'0.45%"
$ws.Range("E14").Style = $origStyle

$ws.Range("B15").Value = "LEO"

$ws.Range("C15").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"

$origStyle = $ws.Range("D15").Style
$ws.Range("D15").Value = "'3.599"
$ws.Range("D15").Style = $origStyle

$origStyle = $ws.Range("E15").Style
$ws.Range("E15").Value = "'2.92%"
$ws.Range("E15").Style = $origStyle

$ws.Range("B16").Value = "GateToken"

$ws.Range("C16").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"

$origStyle = $ws.Range("D16").Style
$ws.Range("D16").Value = "'3.024"
$ws.Range("D16").Style = $origStyle

$origStyle = $ws.Range("E16").Style
$ws.Range("E16").Value = "'0.39%"
$ws.Range("E16").Style = $origStyle

$ws.Range("B17").Value = "BTSEToken"

$ws.Range("C17").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"

$origStyle = $ws.Range("D17").Style
$ws.Range("D17").Value = "'2.055"
$ws.Range("D17").Style = $origStyle

$origStyle = $ws.Range("E17").Style
$ws.Range("E17").Value = "'-1.73%"
$ws.Range("E17").Style = $origStyle

$ws.Range("B18").Value = "BitpandaEcosystemToken"

$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"

$origStyle = $ws.Range("D18").Style
$ws.Range("D18").Value = "'0.3112"
$ws.Range("D18").Style = $origStyle

$origStyle = $ws.Range("E18").Style
$ws.Range("E18").Value = "'-2.40%"
$ws.Range("E18").Style = $origStyle

$ws.Range("B19").Value = "WazirX"

$ws.Range("C19").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"

$origStyle = $ws.Range("D19").Style
$ws.Range("D19").Value = "'0.1330"
$ws.Range("D19").Style = $origStyle

$origStyle = $ws.Range("E19").Style
$ws.Range("E19").Value = "'-0.44%"
$ws.Range("E19").Style = $origStyle

$origStyle = $ws.Range("D20").Style
$ws.Range("D20").Value = "'0.03167"
$ws.Range("D20").Style = $origStyle

$origStyle = $ws.Range("E20").Style
$ws.Range("E20").Value = "'-2.54%"
$ws.Range("E20").Style = $origStyle

$origStyle = $ws.Range("E21").Style
$ws.Range("E21").Value = "'-1.96%"
$ws.Range("E21").Style = $origStyle

$origStyle = $ws.Range("D22").Style
$ws.Range("D22").Value = "'3.744"
$ws.Range("D22").Style = $origStyle

$origStyle = $ws.Range("E22").Style
$ws.Range("E22").Value = "'0.07%"
$ws.Range("E22").Style = $origStyle

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").Value = "'0.04651"
$ws.Range("D23").Style = $origStyle

$origStyle = $ws.Range("E23").Style
$ws.Range("E23").Value = "'-0.93%"
$ws.Range("E23").Style = $origStyle

$origStyle = $ws.Range("D24").Style
$ws.Range("D24").Value = "'0.1374"
$ws.Range("D24").Style = $origStyle

$origStyle = $ws.Range("E24").Style
$ws.Range("E24").Value = "'2.45%"
$ws.Range("E24").Style = $origStyle

$origStyle = $ws.Range("E25").Style
$ws.Range("E25").Value = "'0.30%"
$ws.Range("E25").Style = $origStyle

$origStyle = $ws.Range("D26").Style
$ws.Range("D26").Value = "'0.004546"
$ws.Range("D26").Style = $origStyle

$origStyle = $ws.Range("E26").Style
$ws.Range("E26").Value = "'5.98%"
$ws.Range("E26").Style = $origStyle

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").Value = "'0.00009600"
$ws.Range("D27").Style = $origStyle

$origStyle = $ws.Range("E27").Style
$ws.Range("E27").Value = "'-0.98%"
$ws.Range("E27").Style = $origStyle

$origStyle = $ws.Range("E28").Style
$ws.Range("E28").Value = "'-0.09%"
$ws.Range("E28").Style = $origStyle

$origStyle = $ws.Range("D40").Style
$ws.Range("D40").Value = "'0.03651"
$ws.Range("D40").Style = $origStyle

$origStyle = $ws.Range("E40").Style
$ws.Range("E40").Value = "'-0.39%"
$ws.Range("E40").Style = $origStyle

$origStyle = $ws.Range("D41").Style
$ws.Range("D41").Value = "'0.006150"
$ws.Range("D41").Style = $origStyle

$origStyle = $ws.Range("E41").Style
$ws.Range("E41").Value = "'81.82%"
$ws.Range("E41").Style = $origStyle

$origStyle = $ws.Range("D42").Style
$ws.Range("D42").Value = "'0.1052"
$ws.Range("D42").Style = $origStyle

$origStyle = $ws.Range("E42").Style
$ws.Range("E42").Value = "'-22.11%"
$ws.Range("E42").Style = $origStyle

$origStyle = $ws.Range("D43").Style
$ws.Range("D43").Value = "'0.002596"
$ws.Range("D43").Style = $origStyle

$origStyle = $ws.Range("E43").Style
$ws.Range("E43").Value = "'-4.17%"
$ws.Range("E43").Style = $origStyle

$origStyle = $ws.Range("D44").Style
$ws.Range("D44").Value = "'0.007973"
$ws.Range("D44").Style = $origStyle

$origStyle = $ws.Range("E44").Style
$ws.Range("E44").Value = "'-3.05%"
$ws.Range("E44").Style = $origStyle

$origStyle = $ws.Range("D45").Style
$ws.Range("D45").Value = "'0.00005390"
$ws.Range("D45").Style = $origStyle

$origStyle = $ws.Range("E45").Style
$ws.Range("E45").Value = "'1.84%"
$ws.Range("E45").Style = $origStyle

$origStyle = $ws.Range("D48").Style
$ws.Range("D48").Value = "'0.002402"
$ws.Range("D48").Style = $origStyle

$origStyle = $ws.Range("E48").Style
$ws.Range("E48").Value = "'19.16%"
$ws.Range("E48").Style = $origStyle
